# RBA v2.5 - Atualizacao da Tela
#
# Replaces placeholder text "tre"/"Tre"/"TRE"/"TERE" (various case forms)
# with "qwer"/"Qwer"/"Qewr"/"QWER" equivalents in the document body and
# the primary header, matching each run individually (case-sensitive,
# whole-word) and in document order.

$d = $word.ActiveDocument

# wdReplaceOne = 1, wdFindWrapContinue (no wrap past the search start) = 1
# Find.Execute signature:
#  (FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#   MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)

# --- Body (word/document.xml): single bold "TERE" -> "QWER" ---
$body = $d.Content
$body.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# --- Header (word/header1.xml): sequential occurrences in document order ---
$hdr = $d.Sections(1).Headers(1).Range

$hdr.Find.Execute("TRE",  $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null
$hdr.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null
$hdr.Find.Execute("Tre",  $true, $true, $false, $false, $false, $true, 1, $false, "Qwer", 1) | Out-Null
$hdr.Find.Execute("Tre",  $true, $true, $false, $false, $false, $true, 1, $false, "Qwer", 1) | Out-Null
$hdr.Find.Execute("Tre",  $true, $true, $false, $false, $false, $true, 1, $false, "Qewr", 1) | Out-Null
$hdr.Find.Execute("Tre",  $true, $true, $false, $false, $false, $true, 1, $false, "Qewr", 1) | Out-Null
$hdr.Find.Execute("Tre",  $true, $true, $false, $false, $false, $true, 1, $false, "Qwer", 1) | Out-Null
$hdr.Find.Execute("tre",  $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null
$hdr.Find.Execute("tre",  $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null
$hdr.Find.Execute("tre",  $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null
